$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added at the top of the existing
# "Poroto verde" time series (row 106). Every existing record from row
# 106 down to the former last row (118) shifts down by one row, so we
# insert a fresh row at 106 (carrying the date-format style along) and
# then populate it with the new data point.
$ws.Rows.Item(106).Insert()

$ws.Cells.Item(106, 1).Value = 4
$ws.Cells.Item(106, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(106, 3).Value = "Los Lagos"
$ws.Cells.Item(106, 4).Value = 44946
$ws.Cells.Item(106, 5).Value = 10
$ws.Cells.Item(106, 6).Value = 100112031
$ws.Cells.Item(106, 7).Value = "Poroto verde"
$ws.Cells.Item(106, 8).Value = "Magnum"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 35
$ws.Cells.Item(106, 11).Value = 28000
$ws.Cells.Item(106, 12).Value = 28000
$ws.Cells.Item(106, 13).Value = 28000
$ws.Cells.Item(106, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(106, 15).Value = "Región Metropolitana"
$ws.Cells.Item(106, 16).Value = 1120
$ws.Cells.Item(106, 17).Value = 25
$ws.Cells.Item(106, 18).Value = "Hortaliza"
